$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data with new randomized values, keep Title/Role unchanged
$ws.Range("A2").Value = "TXrGu601"
$ws.Range("B2").Value = 23100515
$ws.Range("C2").Value = "alaehlb21"
$ws.Range("D2").Value = "F35`$xn#H"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "RhsTwfaH"
$ws.Range("G2").Value = "jVCP"
$ws.Range("H2").Value = "Candidate"

# Remove rows 3 and 4 (extra candidate test data no longer needed)
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Refresh the selection to match the new used range
$ws.Range("A1:H2").Select()
